# "XML Item database working"
#
# Fill in the attendance numbers for the "Week 10" sheet (previously all
# blank) and note that Janneke had a doctor's appointment that week, then
# leave the selection/active-sheet state the way the author left it
# (Percentages active, selected at F23; Week 10 selected at F22).

$wb = $excel.ActiveWorkbook

$wsWeek10 = $wb.Worksheets.Item("Week 10")

# Mon (row 2) - only Gerwin (column B) logged hours that day.
$wsWeek10.Range("B2").Value = 0

# Tue (row 3)
$wsWeek10.Range("B3").Value = 4
$wsWeek10.Range("C3").Value = 4
$wsWeek10.Range("D3").Value = 4
$wsWeek10.Range("E3").Value = 0
$wsWeek10.Range("F3").Value = 4
$wsWeek10.Range("G3").Value = 4
$wsWeek10.Range("H3").Value = 4
$wsWeek10.Range("J3").Value = "Janneke had een doktors afspraak."

# Wed (row 4)
$wsWeek10.Range("B4:H4").Value = 2

# Thu (row 5)
$wsWeek10.Range("B5:H5").Value = 4

# Fri (row 6)
$wsWeek10.Range("B6:H6").Value = 8

# Extra "Geoorloofd" (excused) hour note for the week.
$wsWeek10.Range("E12").Value = 4

# Restore the selection on Week 10 (it was the active tab before; the
# author moved the cursor to F22 and switched to the Percentages tab).
$wsWeek10.Range("F22").Select()

$wsPercentages = $wb.Worksheets.Item("Percentages")
$wsPercentages.Activate()
$wsPercentages.Range("F23").Select()
